$d = $word.ActiveDocument

$d.Content.Find.Execute("26+26=52", $true, $false, $false, $false, $false, $true, 1, $false, "48+27=75", 2) | Out-Null
$d.Content.Find.Execute("6+18=24", $true, $false, $false, $false, $false, $true, 1, $false, "16+56=72", 2) | Out-Null
$d.Content.Find.Execute("82-68=14", $true, $false, $false, $false, $false, $true, 1, $false, "47+18=65", 2) | Out-Null
$d.Content.Find.Execute("81-76=5", $true, $false, $false, $false, $false, $true, 1, $false, "58-39=19", 2) | Out-Null
$d.Content.Find.Execute("65+16=81", $true, $false, $false, $false, $false, $true, 1, $false, "31-13=18", 2) | Out-Null
$d.Content.Find.Execute("33-19=14", $true, $false, $false, $false, $false, $true, 1, $false, "40-9=31", 2) | Out-Null
$d.Content.Find.Execute("61-16=45", $true, $false, $false, $false, $false, $true, 1, $false, "78-9=69", 2) | Out-Null
$d.Content.Find.Execute("58+17=75", $true, $false, $false, $false, $false, $true, 1, $false, "92-48=44", 2) | Out-Null
$d.Content.Find.Execute("52+19=71", $true, $false, $false, $false, $false, $true, 1, $false, "28+7=35", 2) | Out-Null
$d.Content.Find.Execute("37+17=54", $true, $false, $false, $false, $false, $true, 1, $false, "34+49=83", 2) | Out-Null
$d.Content.Find.Execute("70-5=65", $true, $false, $false, $false, $false, $true, 1, $false, "92-33=59", 2) | Out-Null
$d.Content.Find.Execute("87+5=92", $true, $false, $false, $false, $false, $true, 1, $false, "17+64=81", 2) | Out-Null
$d.Content.Find.Execute("50-9=41", $true, $false, $false, $false, $false, $true, 1, $false, "58+5=63", 2) | Out-Null
$d.Content.Find.Execute("38-19=19", $true, $false, $false, $false, $false, $true, 1, $false, "68-39=29", 2) | Out-Null
$d.Content.Find.Execute("76-67=9", $true, $false, $false, $false, $false, $true, 1, $false, "72-68=4", 2) | Out-Null
$d.Content.Find.Execute("18+9=27", $true, $false, $false, $false, $false, $true, 1, $false, "70-63=7", 2) | Out-Null
$d.Content.Find.Execute("47+6=53", $true, $false, $false, $false, $false, $true, 1, $false, "92-16=76", 2) | Out-Null
$d.Content.Find.Execute("31-26=5", $true, $false, $false, $false, $false, $true, 1, $false, "7+48=55", 2) | Out-Null
$d.Content.Find.Execute("70-49=21", $true, $false, $false, $false, $false, $true, 1, $false, "75-39=36", 2) | Out-Null
$d.Content.Find.Execute("83-77=6", $true, $false, $false, $false, $false, $true, 1, $false, "39+4=43", 2) | Out-Null
$d.Content.Find.Execute("90-86=4", $true, $false, $false, $false, $false, $true, 1, $false, "80-63=17", 2) | Out-Null
$d.Content.Find.Execute("5+27=32", $true, $false, $false, $false, $false, $true, 1, $false, "60-44=16", 2) | Out-Null
$d.Content.Find.Execute("84-68=16", $true, $false, $false, $false, $false, $true, 1, $false, "14+49=63", 2) | Out-Null
$d.Content.Find.Execute("29+67=96", $true, $false, $false, $false, $false, $true, 1, $false, "28+4=32", 2) | Out-Null
$d.Content.Find.Execute("61-57=4", $true, $false, $false, $false, $false, $true, 1, $false, "19+7=26", 2) | Out-Null
$d.Content.Find.Execute("62-36=26", $true, $false, $false, $false, $false, $true, 1, $false, "54-8=46", 2) | Out-Null
$d.Content.Find.Execute("29+2=31", $true, $false, $false, $false, $false, $true, 1, $false, "46+15=61", 2) | Out-Null
$d.Content.Find.Execute("63-55=8", $true, $false, $false, $false, $false, $true, 1, $false, "85-78=7", 2) | Out-Null
$d.Content.Find.Execute("48+48=96", $true, $false, $false, $false, $false, $true, 1, $false, "39+18=57", 2) | Out-Null
$d.Content.Find.Execute("28+48=76", $true, $false, $false, $false, $false, $true, 1, $false, "42-34=8", 2) | Out-Null
$d.Content.Find.Execute("67-19=48", $true, $false, $false, $false, $false, $true, 1, $false, "35+28=63", 2) | Out-Null
$d.Content.Find.Execute("17+78=95", $true, $false, $false, $false, $false, $true, 1, $false, "81-39=42", 2) | Out-Null
$d.Content.Find.Execute("78+14=92", $true, $false, $false, $false, $false, $true, 1, $false, "18+74=92", 2) | Out-Null
$d.Content.Find.Execute("69+29=98", $true, $false, $false, $false, $false, $true, 1, $false, "92-3=89", 2) | Out-Null
$d.Content.Find.Execute("8+14=22", $true, $false, $false, $false, $false, $true, 1, $false, "21-19=2", 2) | Out-Null
$d.Content.Find.Execute("47-39=8", $true, $false, $false, $false, $false, $true, 1, $false, "90-69=21", 2) | Out-Null
$d.Content.Find.Execute("52-18=34", $true, $false, $false, $false, $false, $true, 1, $false, "25+36=61", 2) | Out-Null
$d.Content.Find.Execute("71-38=33", $true, $false, $false, $false, $false, $true, 1, $false, "54-25=29", 2) | Out-Null
$d.Content.Find.Execute("46+39=85", $true, $false, $false, $false, $false, $true, 1, $false, "73-58=15", 2) | Out-Null
$d.Content.Find.Execute("14+48=62", $true, $false, $false, $false, $false, $true, 1, $false, "84+7=91", 2) | Out-Null
$d.Content.Find.Execute("68+6=74", $true, $false, $false, $false, $false, $true, 1, $false, "8+48=56", 2) | Out-Null
$d.Content.Find.Execute("58+4=62", $true, $false, $false, $false, $false, $true, 1, $false, "7+48=55", 2) | Out-Null
$d.Content.Find.Execute("70-38=32", $true, $false, $false, $false, $false, $true, 1, $false, "33-17=16", 2) | Out-Null
$d.Content.Find.Execute("9+4=13", $true, $false, $false, $false, $false, $true, 1, $false, "72-56=16", 2) | Out-Null
$d.Content.Find.Execute("39+25=64", $true, $false, $false, $false, $false, $true, 1, $false, "90-87=3", 2) | Out-Null
$d.Content.Find.Execute("68+8=76", $true, $false, $false, $false, $false, $true, 1, $false, "8+49=57", 2) | Out-Null
$d.Content.Find.Execute("62-46=16", $true, $false, $false, $false, $false, $true, 1, $false, "25-19=6", 2) | Out-Null
$d.Content.Find.Execute("45-19=26", $true, $false, $false, $false, $false, $true, 1, $false, "44+7=51", 2) | Out-Null
$d.Content.Find.Execute("91-66=25", $true, $false, $false, $false, $false, $true, 1, $false, "85-68=17", 2) | Out-Null
$d.Content.Find.Execute("40-21=19", $true, $false, $false, $false, $false, $true, 1, $false, "82-79=3", 2) | Out-Null
$d.Content.Find.Execute("7+49=56", $true, $false, $false, $false, $false, $true, 1, $false, "27+37=64", 2) | Out-Null
$d.Content.Find.Execute("47+34=81", $true, $false, $false, $false, $false, $true, 1, $false, "92-29=63", 2) | Out-Null
$d.Content.Find.Execute("35+29=64", $true, $false, $false, $false, $false, $true, 1, $false, "82-78=4", 2) | Out-Null
$d.Content.Find.Execute("81-12=69", $true, $false, $false, $false, $false, $true, 1, $false, "90-45=45", 2) | Out-Null
$d.Content.Find.Execute("36+58=94", $true, $false, $false, $false, $false, $true, 1, $false, "82-79=3", 2) | Out-Null
$d.Content.Find.Execute("25+7=32", $true, $false, $false, $false, $false, $true, 1, $false, "49+45=94", 2) | Out-Null
$d.Content.Find.Execute("18+63=81", $true, $false, $false, $false, $false, $true, 1, $false, "53-24=29", 2) | Out-Null
$d.Content.Find.Execute("17+27=44", $true, $false, $false, $false, $false, $true, 1, $false, "37+56=93", 2) | Out-Null
$d.Content.Find.Execute("85-17=68", $true, $false, $false, $false, $false, $true, 1, $false, "4+8=12", 2) | Out-Null
$d.Content.Find.Execute("13+8=21", $true, $false, $false, $false, $false, $true, 1, $false, "18+26=44", 2) | Out-Null
$d.Content.Find.Execute("18+23=41", $true, $false, $false, $false, $false, $true, 1, $false, "5+57=62", 2) | Out-Null
$d.Content.Find.Execute("22-9=13", $true, $false, $false, $false, $false, $true, 1, $false, "92-9=83", 2) | Out-Null
$d.Content.Find.Execute("63-28=35", $true, $false, $false, $false, $false, $true, 1, $false, "22-17=5", 2) | Out-Null
$d.Content.Find.Execute("17+18=35", $true, $false, $false, $false, $false, $true, 1, $false, "40-27=13", 2) | Out-Null
$d.Content.Find.Execute("49+22=71", $true, $false, $false, $false, $false, $true, 1, $false, "87-39=48", 2) | Out-Null
$d.Content.Find.Execute("57+26=83", $true, $false, $false, $false, $false, $true, 1, $false, "27+67=94", 2) | Out-Null
$d.Content.Find.Execute("44+48=92", $true, $false, $false, $false, $false, $true, 1, $false, "18+8=26", 2) | Out-Null
$d.Content.Find.Execute("36+38=74", $true, $false, $false, $false, $false, $true, 1, $false, "6+59=65", 2) | Out-Null
$d.Content.Find.Execute("16+48=64", $true, $false, $false, $false, $false, $true, 1, $false, "38+15=53", 2) | Out-Null
$d.Content.Find.Execute("80-58=22", $true, $false, $false, $false, $false, $true, 1, $false, "89+6=95", 2) | Out-Null
$d.Content.Find.Execute("19+25=44", $true, $false, $false, $false, $false, $true, 1, $false, "7+65=72", 2) | Out-Null
$d.Content.Find.Execute("50-33=17", $true, $false, $false, $false, $false, $true, 1, $false, "46-28=18", 2) | Out-Null
$d.Content.Find.Execute("7+85=92", $true, $false, $false, $false, $false, $true, 1, $false, "73-5=68", 2) | Out-Null
$d.Content.Find.Execute("14+69=83", $true, $false, $false, $false, $false, $true, 1, $false, "48+14=62", 2) | Out-Null
$d.Content.Find.Execute("40-38=2", $true, $false, $false, $false, $false, $true, 1, $false, "53-25=28", 2) | Out-Null
$d.Content.Find.Execute("4+89=93", $true, $false, $false, $false, $false, $true, 1, $false, "68-9=59", 2) | Out-Null
$d.Content.Find.Execute("52-39=13", $true, $false, $false, $false, $false, $true, 1, $false, "44-5=39", 2) | Out-Null
$d.Content.Find.Execute("72-5=67", $true, $false, $false, $false, $false, $true, 1, $false, "37+29=66", 2) | Out-Null
$d.Content.Find.Execute("19+17=36", $true, $false, $false, $false, $false, $true, 1, $false, "71-42=29", 2) | Out-Null
$d.Content.Find.Execute("19+53=72", $true, $false, $false, $false, $false, $true, 1, $false, "29+54=83", 2) | Out-Null
$d.Content.Find.Execute("9+38=47", $true, $false, $false, $false, $false, $true, 1, $false, "15+17=32", 2) | Out-Null
$d.Content.Find.Execute("94-66=28", $true, $false, $false, $false, $false, $true, 1, $false, "26+35=61", 2) | Out-Null
$d.Content.Find.Execute("19+29=48", $true, $false, $false, $false, $false, $true, 1, $false, "98-49=49", 2) | Out-Null
$d.Content.Find.Execute("95-28=67", $true, $false, $false, $false, $false, $true, 1, $false, "19+37=56", 2) | Out-Null
$d.Content.Find.Execute("14-5=9", $true, $false, $false, $false, $false, $true, 1, $false, "50-36=14", 2) | Out-Null
$d.Content.Find.Execute("62-7=55", $true, $false, $false, $false, $false, $true, 1, $false, "15+59=74", 2) | Out-Null
$d.Content.Find.Execute("56+28=84", $true, $false, $false, $false, $false, $true, 1, $false, "61-28=33", 2) | Out-Null
$d.Content.Find.Execute("34+28=62", $true, $false, $false, $false, $false, $true, 1, $false, "51-12=39", 2) | Out-Null
$d.Content.Find.Execute("9+77=86", $true, $false, $false, $false, $false, $true, 1, $false, "61-15=46", 2) | Out-Null
$d.Content.Find.Execute("77+6=83", $true, $false, $false, $false, $false, $true, 1, $false, "37+26=63", 2) | Out-Null
$d.Content.Find.Execute("81-52=29", $true, $false, $false, $false, $false, $true, 1, $false, "52-45=7", 2) | Out-Null
$d.Content.Find.Execute("83-16=67", $true, $false, $false, $false, $false, $true, 1, $false, "39+52=91", 2) | Out-Null
$d.Content.Find.Execute("80-66=14", $true, $false, $false, $false, $false, $true, 1, $false, "49+13=62", 2) | Out-Null
$d.Content.Find.Execute("9+18=27", $true, $false, $false, $false, $false, $true, 1, $false, "51-3=48", 2) | Out-Null
$d.Content.Find.Execute("94-49=45", $true, $false, $false, $false, $false, $true, 1, $false, "96-78=18", 2) | Out-Null
$d.Content.Find.Execute("72+19=91", $true, $false, $false, $false, $false, $true, 1, $false, "78+3=81", 2) | Out-Null
$d.Content.Find.Execute("92-13=79", $true, $false, $false, $false, $false, $true, 1, $false, "35+6=41", 2) | Out-Null
$d.Content.Find.Execute("48+16=64", $true, $false, $false, $false, $false, $true, 1, $false, "14+57=71", 2) | Out-Null
$d.Content.Find.Execute("46+36=82", $true, $false, $false, $false, $false, $true, 1, $false, "94-69=25", 2) | Out-Null
$d.Content.Find.Execute("23-5=18", $true, $false, $false, $false, $false, $true, 1, $false, "38+58=96", 2) | Out-Null
